# ---------------------------------------------------------------------------
# Applies the "Produccion 5 / Universal" revision to
# "Compatibilidad Producciones.xlsx":
#   - HARDWARE: bump the revision counter, add a 5th "Produccion 5" /
#     "Universal" block (columns N:P) mirroring the existing Produccion 4
#     block, insert a new "Universal" antenna row, and fill in the newly
#     required "-" / footnote cells in the bottom notes row.
#   - SOFTWARE: bump the version/date counters and the three running IDs.
#   - Make HARDWARE the active sheet/tab, matching the saved selection.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$hw  = $wb.Worksheets.Item("HARDWARE")
$sw  = $wb.Worksheets.Item("SOFTWARE")

# ---------------------------------------------------------------------------
# HARDWARE sheet
# ---------------------------------------------------------------------------

# Revision counter on B2 (3 -> 4); B3 keeps its =TODAY() formula as-is.
$hw.Range("B2").Value = 4

# --- New "Produccion 5" block (columns N:P), cloned from the "Produccion 4"
# block (columns K:M) so the fill/border/font formatting matches exactly. ---
$hw.Range("K5:M9").Copy()
$hw.Range("N5").PasteSpecial(-4122)

$hw.Range("N5:P5").Merge()
$hw.Range("N5").Value = "Produccion 5"

$hw.Range("N6").Value = "CPU"
$hw.Range("O6").Value = "DCDC"
$hw.Range("P6").Value = "CONECTORES"

$hw.Range("O9").Value = "*"

# --- Insert the new "Universal" antenna row above the notes row. ---
$hw.Rows.Item(10).Insert()

$hw.Range("A8:M8").Copy()
$hw.Range("A10").PasteSpecial(-4122)
$hw.Rows.Item(10).RowHeight = 21.75
$hw.Range("A10").Value = "Universal"

# Clone the Produccion-4 formatting (now K10:M11) into the new Produccion-5
# columns for the just-inserted row 10 and the shifted notes row 11.
$hw.Range("K10:M11").Copy()
$hw.Range("N10").PasteSpecial(-4122)

# --- Fill in the newly required cells on the (shifted) notes row 11. ---
$hw.Range("I11").Value = "-"
$hw.Range("J11").Value = "-"
$hw.Range("M11").Value = "-"
$hw.Range("N11").Value = "-"
$hw.Range("O11").Value = " * Res. Balastro`nR19,R20 = 20R"

$hw.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# SOFTWARE sheet
# ---------------------------------------------------------------------------

$sw.Range("B2").Value = 2
$sw.Range("B3").Value = 45547

$sw.Range("A7").Value = 125
$sw.Range("B7").Value = 313
$sw.Range("D7").Value = 205

# ---------------------------------------------------------------------------
# Active sheet / selection bookkeeping
# ---------------------------------------------------------------------------

$sw.Range("C18").Select()
$hw.Activate()
$hw.Range("M13").Select()
